# Update the residue-combination / distance table with new charge-based
# distances and frame counts, then remove the now-obsolete trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the replacement values in columns C/D (and occasionally A) are
# plain digit strings (e.g. "1105", "455", "6017"). Excel would normally
# auto-convert those to numbers when assigned through .Value, but the
# source data keeps them as text, so force a text number format on those
# specific cells before writing the values.
$textCells = @("C9", "C13", "D13", "C14", "D14", "C15", "D15", "D16", "D17", "C18", "D18")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Update existing rows with new values -------------------------------

# Row 2
$ws.Range("A2").Value = "130, 130, 780, 1105"

# Row 3
$ws.Range("A3").Value = "130, 455, 780, 780"

# Row 5
$ws.Range("A5").Value = "130, 130, 780, 780"

# Row 7
$ws.Range("A7").Value = "455, 780, 1105, 1105"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = "1105, 1105, 1105, 1105"
$ws.Range("D7").Value = "5131, 5331, 5433, 6562"

# Row 8
$ws.Range("A8").Value = "130, 130, 455, 1105"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "130, 130, 130"
$ws.Range("D8").Value = "4416, 6489, 6670"

# Row 9
$ws.Range("A9").Value = "455, 780, 1105, 1105, SF"
$ws.Range("C9").Value = "1105"

# Row 10
$ws.Range("A10").Value = "130, 455, 780, 1105"

# Row 11
$ws.Range("A11").Value = "130, 780, 1105, 1105"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "1105, 1105"
$ws.Range("D11").Value = "5400, 6359"

# Row 12
$ws.Range("A12").Value = "455, 780, 780, 1105"

# Row 13
$ws.Range("A13").Value = "455, 455, 780, 1105"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "455"
$ws.Range("D13").Value = "5886"

# Row 14
$ws.Range("A14").Value = "130, 130, 780, 780, 1105"
$ws.Range("C14").Value = "130"
$ws.Range("D14").Value = "6017"

# Row 15
$ws.Range("A15").Value = "455, 780, 780, 1105, SF"
$ws.Range("C15").Value = "780"
$ws.Range("D15").Value = "6202"

# Row 16
$ws.Range("A16").Value = "130, 130, 455"
$ws.Range("D16").Value = "6427"

# Row 17
$ws.Range("A17").Value = "130, 130, 130, 455, 780"
$ws.Range("D17").Value = "6727"

# Row 18
$ws.Range("A18").Value = "130, 455, 780, SF"
$ws.Range("C18").Value = "780"
$ws.Range("D18").Value = "6748"

# --- Remove obsolete trailing rows 19-22 --------------------------------

$ws.Range("A19:D22").Delete()
